$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.644.59"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.487.27"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "2.949.87"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "67.570.08"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "2.475.40"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "2.615.14"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "506.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +6.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0259"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0745"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.586"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
